# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 / Row 17: swap the "Periodo Mora" labels (2504 <-> 2503) and
# update "Salario Basico" (F) / "Valor Mora" (G) values accordingly.
$ws.Range("E16").Value = "2503"
$ws.Range("F16").Value = 11388
$ws.Range("G16").Value = 1540000

$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1540000
